$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the location keyword value in B2 from "San+Diego" to "San+Franciso"
$ws.Range("B2").Value = "San+Franciso"

# Reflect the new active selection on Sheet1 (user ended editing on B2)
$ws.Activate()
$ws.Range("B2").Select()
